# Work diary update for the day: new commit id for the GanttProject remark,
# a Scrum-sandbox entry and a documentation-canvas entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "JdT-TPI_LRD"

# Row 8: the GanttProject remark gets its own commit id in column F (kept
# in the same order the rest of the table uses: commit id then remark), and
# the row grows taller to fit the longer, wrapped remark text.
$ws.Range("F8").Value = "9b11e557a8f607254b2b476a41fe98b5022022df"

# Row 9: new entry about filling the Scrum sandbox, with a link to the
# backlog.
$ws.Range("D9").Value = "Remplissage de la sandbox sur Scrum"
$ws.Range("E9").Value = "https://icescrum.cpnv.ch/p/NEWSWEBSIT/"

# Row 10: new entry about adding the documentation canvas, with its commit.
$ws.Range("D10").Value = "Ajout du canvas de documentation"
$ws.Range("F10").Value = "3b50924c0668f0b7c5f34f4eb16db6148d4fd43b"

# Row 8's remark text is rewritten last so the new shared strings land in
# the same allocation order as the rest of the day's entries.
$ws.Range("E8").Value = "Sur GanttProject. Je ne suis pas 100% sûr du resultat. Elle sera révisée au besoin."
$ws.Rows("8").RowHeight = 30

# Dates for the two new rows - copy the date cell's format (A8) so they
# share the existing date style instead of minting a new one.
$ws.Range("A8").Copy()
$ws.Range("A9:A10").PasteSpecial(-4122)
$ws.Range("A9").Value = 44683
$ws.Range("A10").Value = 44683

$ws.Range("C9").Value = 0.75
$ws.Range("C10").Value = 0

# Hyperlink the Scrum backlog URL, reusing the same format as the existing
# hyperlink cell (E3) so it doesn't mint a new style.
$ws.Hyperlinks.Add($ws.Range("E9"), "https://icescrum.cpnv.ch/p/NEWSWEBSIT/")
$ws.Range("E3").Copy()
$ws.Range("E9").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

[void]$ws.Range("E8").Select()
